$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.0984234321627112
$ws.Range("B3").Value = 0.1047075183849701
$ws.Range("H3").Value = 0.2031309505476813
$ws.Range("B4").Value = 0.09897717964829653
$ws.Range("H4").Value = 0.1974006118110077
$ws.Range("B5").Value = 0.07947544163689636
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = 0.1778988737996076
$ws.Range("B6").Value = 0.05182823082317872
$ws.Range("C6").Value = 0.003393416648242191
$ws.Range("D6").Value = 10.33050640383841
$ws.Range("E6").Value = 0.01441222100071589
$ws.Range("F6").Value = 0.04516837039422728
$ws.Range("G6").Value = 0.05848809125213079
$ws.Range("H6").Value = 0.1502516629858899
$ws.Range("B7").Value = 0.05215327718220823
$ws.Range("H7").Value = 0.1505767093449194
$ws.Range("B8").Value = 0.0505871872593374
$ws.Range("C8").Value = 0.003243509853448842
$ws.Range("D8").Value = 5.580092923639919
$ws.Range("E8").Value = 0.01785842071342533
$ws.Range("F8").Value = 0.04422007702240081
$ws.Range("G8").Value = 0.05695429749627432
$ws.Range("H8").Value = 0.1490106194220486
$ws.Range("B9").Value = 0.05372497167157138
$ws.Range("C9").Value = 0.002811005531344026
$ws.Range("D9").Value = 4.618073889692631
$ws.Range("E9").Value = 0.01395161275276387
$ws.Range("F9").Value = 0.04821209580664321
$ws.Range("G9").Value = 0.05923784753650124
$ws.Range("H9").Value = 0.1521484038342826
$ws.Range("B10").Value = 0.0558167460178016
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 0.1542401781805128
$ws.Range("B11").Value = 0.02845956604289967
$ws.Range("H11").Value = 0.1268829982056109
$ws.Range("B12").Value = 0.04965314735332428
$ws.Range("H12").Value = 0.1480765795160355
$ws.Range("B13").Value = 0.06373438030096099
$ws.Range("H13").Value = 0.1621578124636722
$ws.Range("B14").Value = 0.07104767068378057
$ws.Range("H14").Value = 0.1694711028464918
$ws.Range("B15").Value = 0.07553372761543573
$ws.Range("H15").Value = 0.1739571597781469
$ws.Range("B16").Value = 0.0799038185188436
$ws.Range("H16").Value = 0.1783272506815548
$ws.Range("B17").Value = 0.08282733730190932
$ws.Range("H17").Value = 0.1812507694646205
$ws.Range("B18").Value = -0.0984234321627112
$ws.Range("C18").Value = 0.01068438068956193
$ws.Range("D18").Value = -17.27768835861711
$ws.Range("E18").Value = 0.04626682776618744
$ws.Range("F18").Value = -0.1194315750007144
$ws.Range("G18").Value = -0.07741528932470813
$ws.Range("H18").Value = 0
$ws.Range("B19").Value = 0.08587118883494811
$ws.Range("H19").Value = 0.1842946209976593
$ws.Range("B20").Value = 0.08633252247960442
$ws.Range("H20").Value = 0.1847559546423156
$ws.Range("B21").Value = 0.0927173209213141
$ws.Range("H21").Value = 0.1911407530840253
$ws.Range("B22").Value = 0.095347666107395
$ws.Range("H22").Value = 0.1937710982701062
$ws.Range("B23").Value = 0.09838714458859858
$ws.Range("H23").Value = 0.1968105767513098
$ws.Range("B24").Value = 0.1013695941648085
$ws.Range("C24").Value = 0.008373327976158184
$ws.Range("D24").Value = 20.72497380922693
$ws.Range("E24").Value = 0.05754622721266231
$ws.Range("F24").Value = 0.0849002227771836
$ws.Range("G24").Value = 0.117838965552433
$ws.Range("H24").Value = 0.1997930263275197
$ws.Range("B25").Value = 0.101888784563982
$ws.Range("C25").Value = 0.008371035866017383
$ws.Range("D25").Value = 20.46220124530411
$ws.Range("E25").Value = 0.05954346281266894
$ws.Range("F25").Value = 0.08542163044845742
$ws.Range("G25").Value = 0.1183559386795065
$ws.Range("H25").Value = 0.2003122167266932
$ws.Range("B26").Value = 0.1035065334072906
$ws.Range("C26").Value = 0.008286485437668741
$ws.Range("D26").Value = 20.26555869680995
$ws.Range("E26").Value = 0.06571045960041402
$ws.Range("F26").Value = 0.0872047672537817
$ws.Range("G26").Value = 0.1198082995607995
$ws.Range("H26").Value = 0.2019299655700018
$ws.Range("B27").Value = 0.110136090607751
$ws.Range("C27").Value = 0.008484351742612996
$ws.Range("D27").Value = 19.94007017008146
$ws.Range("E27").Value = 0.0676018973312657
$ws.Range("F27").Value = 0.09344886607715885
$ws.Range("G27").Value = 0.1268233151383435
$ws.Range("H27").Value = 0.2085595227704622
$ws.Range("B28").Value = 0.1107059510211285
$ws.Range("C28").Value = 0.0082153362048246
$ws.Range("D28").Value = 19.91213547012392
$ws.Range("E28").Value = 0.08362679391497141
$ws.Range("F28").Value = 0.0945506648204571
$ws.Range("G28").Value = 0.1268612372218009
$ws.Range("H28").Value = 0.2091293831838397
$ws.Range("B29").Value = 0.05551354628186445
$ws.Range("C29").Value = 0.00367589659326651
$ws.Range("D29").Value = 4.452279470895774
$ws.Range("E29").Value = 0.0058356924756495
$ws.Range("F29").Value = 0.04829476130608929
$ws.Range("G29").Value = 0.06273233125763926
$ws.Range("H29").Value = 0.1539369784445757
